# Build the "readme", "DosePerBodyweight", "DosePerSurfaceArea" sheet layout
$wb = $excel.ActiveWorkbook

$doseWeight = $wb.Worksheets.Item("DosePerBodyweight")

# ---- 1) insert the new "readme" sheet in front of DosePerBodyweight ----
$readme = $wb.Worksheets.Add($doseWeight)
$readme.Name = "readme"

$readme.Range("A1").Value = "First column function handel wich is used to set an application parameter`navailable are @addDosetablePerWeight and @addDosetablePerBSA"
$readme.Range("A2").Value = "second hadle list of parameters which are set by the function"
$readme.Range("A3").Value = "Third line header for numeric info, `nfor addDosetablePerWeight folllwoing columns are mandatory column BWmin, BWmax and targetParameter`nfor @addDosetablePerBSA folllwoing columns are mandatory column BSAmin, BSAmax and targetParameter`n"
$readme.Range("A4").Value = "additional colmuns can be added. Please make sure column name should no contain specialletters, also no spaces"
$readme.Range("A6").Value = "attention MoBi internla Units are used , Body weight [kk], BSA [dm^2], DrugMass µmol "

$readme.Range("A1").WrapText = $true
$readme.Range("A3").WrapText = $true

$readme.Rows.Item(1).RowHeight = 27.6
$readme.Rows.Item(3).RowHeight = 82.8

$readme.Columns.Item(1).ColumnWidth = 71.48493303571429

$readme.Activate() | Out-Null
$readme.Range("A15").Select() | Out-Null

# ---- 2) re-fetch DosePerBodyweight (stale handle after sheet insert) and add "DosePerSurfaceArea" after it ----
$doseWeight = $wb.Worksheets.Item("DosePerBodyweight")
$doseBSA = $wb.Worksheets.Add($null, $doseWeight)
$doseBSA.Name = "DosePerSurfaceArea"

$doseBSA.Range("A1").Value = "functionHandle = @addDosetablePerBSA"
$doseBSA.Range("A2").Value = "targetParameterList = {'*Application_*|ProtocolSchemaItem|DrugMass'}"
$doseBSA.Range("A3").Value = "BSAmin"
$doseBSA.Range("B3").Value = "BSAmax"
$doseBSA.Range("C3").Value = "targetParameter"
$doseBSA.Range("D3").Value = "dose_mg"

$doseBSA.Range("A4").Value = 50
$doseBSA.Range("B4").Value = 60
$doseBSA.Range("D4").Value = 55
$doseBSA.Range("C4").Formula = "=D4/225.21*1000"

$doseBSA.Range("A5").Formula = "=A4+10"
$doseBSA.Range("B5").Formula = "=B4+10"
$doseBSA.Range("D5").Formula = "=D4+10"

# matches the shared-formula grouping used on DosePerBodyweight / DosePerSurfaceArea:
#   C5:C10 -> one shared group, A6:B10 -> one shared group, D6:D10 -> one shared group
$doseBSA.Range("C5:C10").Formula = "=D5/225.21*1000"
$doseBSA.Range("A6:B10").Formula = "=A5+10"
$doseBSA.Range("D6:D10").Formula = "=D5+10"

$doseBSA.Columns.Item(1).ColumnWidth = 13.285714285714286
$doseBSA.Columns.Item(3).ColumnWidth = 16.883370535714285

$doseBSA.PageSetup.PaperSize = 9
$doseBSA.PageSetup.Orientation = 1

$doseBSA.Activate() | Out-Null
$doseBSA.Range("F11").Select() | Out-Null

# ---- 3) final active tab must be "readme" (tabSelected="1" in the diff), cursor on A15 ----
$readme = $wb.Worksheets.Item("readme")
$readme.Activate() | Out-Null
$readme.Range("A15").Select() | Out-Null

Write-Host "done"
